# Apply the "Updated with CB, TVV, MT, GP results" edit.
# Adds four new Superclasico race results (Coppa Bernocchi, Tre Valli
# Varesine, Milano Torino, Gran Piemonte) to both the "Overview" and "UTs"
# worksheets, inserted right above the existing "Classics Squad" block.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("UTs")

# ---------------------------------------------------------------------
# Insert the new blank rows first (formatting is inherited from the row
# immediately above, matching native Excel "Insert" behaviour).
# ---------------------------------------------------------------------
$ws1.Rows.Item(58).Insert() | Out-Null
$ws1.Rows.Item(58).Insert() | Out-Null
$ws1.Rows.Item(58).Insert() | Out-Null
$ws1.Rows.Item(58).Insert() | Out-Null

$ws2.Rows.Item(60).Insert() | Out-Null
$ws2.Rows.Item(60).Insert() | Out-Null
$ws2.Rows.Item(60).Insert() | Out-Null
$ws2.Rows.Item(60).Insert() | Out-Null

# ---------------------------------------------------------------------
# Fill in the "Type" / count columns that only reuse pre-existing text
# (kept early, purely cosmetic - order amongst these does not matter).
# ---------------------------------------------------------------------
$ws1.Cells.Item(58, 1).Value = "Superclasico"
$ws1.Cells.Item(59, 1).Value = "Superclasico"
$ws1.Cells.Item(60, 1).Value = "Superclasico"
$ws1.Cells.Item(61, 1).Value = "Superclasico"

$ws1.Cells.Item(58, 3).Value = 6
$ws1.Cells.Item(58, 4).Value = 34
$ws1.Cells.Item(58, 5).Value = 34
$ws1.Cells.Item(59, 3).Value = 6
$ws1.Cells.Item(59, 4).Value = 40
$ws1.Cells.Item(59, 5).Value = 40
$ws1.Cells.Item(60, 3).Value = 6
$ws1.Cells.Item(60, 4).Value = 46
$ws1.Cells.Item(60, 5).Value = 46
$ws1.Cells.Item(61, 3).Value = 6
$ws1.Cells.Item(61, 4).Value = 45
$ws1.Cells.Item(61, 5).Value = 45

$ws1.Cells.Item(59, 6).Value = "10''"
$ws1.Cells.Item(60, 6).Value = "24''"
$ws1.Cells.Item(61, 6).Value = "21''"

# ---------------------------------------------------------------------
# The sequence below reproduces the exact order in which brand-new
# shared strings were first introduced in the authored workbook.
# ---------------------------------------------------------------------
$ws1.Cells.Item(58, 2).Value = "Coppa Bernocchi"
$ws1.Cells.Item(59, 2).Value = "Tre Valli Varesine"
$ws1.Cells.Item(60, 2).Value = "Milano Torino"
$ws1.Cells.Item(58, 6).Value = "3''"

$ws2.Cells.Item(60, 1).Value = "Superclasico"
$ws2.Cells.Item(60, 2).Value = "Coppa Bernocchi"
$ws2.Cells.Item(60, 3).Value = 601
$ws2.Cells.Item(60, 4).Value = 94
$ws2.Cells.Item(60, 5).Value = "Remco Evenepoel"
$ws2.Cells.Item(60, 6).Value = "Thibaut Pinot"
$ws2.Cells.Item(60, 7).Value = "Fausto Masnada"
$ws2.Cells.Item(60, 8).Value = "Alessandro Covi"
$ws2.Cells.Item(60, 9).Value = "Samuele Battistella"
$ws2.Cells.Item(60, 10).Value = "Juan Sebastian Molano"

$ws2.Cells.Item(61, 1).Value = "Superclasico"
$ws2.Cells.Item(61, 2).Value = "Tre Valli Varesine"
$ws2.Cells.Item(61, 3).Value = 616
$ws2.Cells.Item(61, 4).Value = 86
$ws2.Cells.Item(61, 5).Value = "Tadej Pogacar"
$ws2.Cells.Item(61, 6).Value = "Davide Formolo"
$ws2.Cells.Item(61, 7).Value = "Alessandro De Marchi"
$ws2.Cells.Item(61, 8).Value = "Lorenzo Rota"
$ws2.Cells.Item(61, 9).Value = "Andreas Kron"
$ws2.Cells.Item(61, 10).Value = "Nelson Oliveira"

$ws1.Cells.Item(61, 2).Value = "Gran Piemonte"

$ws2.Cells.Item(62, 1).Value = "Superclasico"
$ws2.Cells.Item(62, 2).Value = "Milano Torino"
$ws2.Cells.Item(62, 3).Value = 416
$ws2.Cells.Item(62, 4).Value = 100
$ws2.Cells.Item(62, 5).Value = "Primoz Roglic"
$ws2.Cells.Item(62, 6).Value = "Adam Yates"
$ws2.Cells.Item(62, 7).Value = "David Gaudu"
$ws2.Cells.Item(62, 8).Value = "Fausto Masnada"
$ws2.Cells.Item(62, 9).Value = "Joao Almeida"
$ws2.Cells.Item(62, 10).Value = "Mauri Vansevenant"

$ws2.Cells.Item(63, 1).Value = "Superclasico"
$ws2.Cells.Item(63, 2).Value = "Gran Piemonte"
$ws2.Cells.Item(63, 3).Value = 455
$ws2.Cells.Item(63, 4).Value = 70
$ws2.Cells.Item(63, 5).Value = "Matteo Trentin"
$ws2.Cells.Item(63, 6).Value = "Giacomo Nizzolo"
$ws2.Cells.Item(63, 7).Value = "Jakub Marezcko"
$ws2.Cells.Item(63, 8).Value = "Olav Kooij"
$ws2.Cells.Item(63, 9).Value = "Biniam Ghirmay"
$ws2.Cells.Item(63, 10).Value = "Matthew Halls"

# ---------------------------------------------------------------------
# Selection state (matches the saved workbook: UTs selection set first,
# then Overview re-activated/selected so it stays the visible tab).
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C63").Select() | Out-Null

$ws1.Activate()
$ws1.Range("F61").Select() | Out-Null
